# Daily attendance processing - 2026-01-24 11:33:47
# Swap the order of "System" and the recorder's email in the
# "Recorded By" column (G) wherever the value is exactly
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
